# Update the "repaymentstrategy" value on the ProductLoanInput sheet
# from "RBI (India)" to "Overdue/Due Fee/Int,Principal".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Leave the active selection on the edited cell, matching the authored workbook.
$ws.Activate()
$ws.Range("B17").Select()
